$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Unmerge D2:O2 (merge cell is being removed)
$ws.Range("D2:O2").UnMerge()

# 2) Row 2 height changes from 11.25 to 15 (still custom height)
$ws.Rows(2).RowHeight = 15

# 3) D2:O2 keep border/font/wrap but drop horizontal=center alignment (now default/general horizontal)
$ws.Range("D2:O2").HorizontalAlignment = 1   # xlGeneral

# 4) New column Q values
# Q2: empty cell, style like P2 (font1/border1) but without fill applied - visually same as P2
$ws.Range("P2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q2").Value = $null

# Q3: 2020, same style as P3
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q3").Value = 2020

# Q4: 14.5, style like P4 but with fill explicitly applied
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q4").Value = 14.5

# Q5: 13.8, same style as P5
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q5").Value = 13.8

$excel.CutCopyMode = $false

# 5) Update selection to match diff (activeCell P13)
$ws.Range("P13").Select()
